# Add cantrals by cantons
# ---------------------------------------------------------------------------
# The sheet used to have a two-row header (a "(m3/s)/(MW)/(GWh)" row and a
# "Hiver/Ete/Annee" row) above three data rows. It becomes a single-row
# header (idx/idx2/Name/Date Start/Date End/(m3/s)/(MW1)/(MW2)/(GWh) Winter/
# (GWh) Summer/(GWh) Year) directly above the same three data rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Create a style equivalent to the existing header style (Arial 9,
#    General format) but without a redundant applied-number-format flag, to
#    use on the new F1:K1 header cells. We add it as a temporary named
#    style so Excel mints a fresh cellXf, then drop the name again (the
#    cellXf itself survives and stays referenced by the cells using it).
# ---------------------------------------------------------------------------
$tmpStyle = $wb.Styles.Add("TmpHeaderStyle")
$tmpStyle.Font.Name = "Arial"
$tmpStyle.Font.Size = 9

# ---------------------------------------------------------------------------
# 2) Drop the old sub-header row (row 2: "Hiver/Ete/Annee" units row). This
#    shifts the three data rows up from 3/4/5 to 2/3/4.
# ---------------------------------------------------------------------------
$ws.Rows(2).Delete()

# ---------------------------------------------------------------------------
# 3) Rewrite row 1 as a single header row.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$ws.Range("F1:K1").Style = "TmpHeaderStyle"

# A1:E1 keep the plain/default style (no explicit style id), matching the
# target layout.
$ws.Range("A1:E1").Style = "Normal"

# Clean up the temporary named style; the cellXf it minted remains in use by
# F1:K1 so the formatting is preserved.
$wb.Styles.Item("TmpHeaderStyle").Delete()

# ---------------------------------------------------------------------------
# 4) Selection lands on the first data row, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("A2:K2").Select() | Out-Null
